$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "mario lau new data": refreshed 客座率Lau (occupancy-rate) series in row 4 ---
$row4 = [ordered]@{
    "C4" = 69.3
    "D4" = 60.1
    "E4" = 51.4
    "F4" = 52.9
    "G4" = 55.1
    "H4" = 56.8
    "I4" = 58.7
    "J4" = 62.5
    "K4" = 65.4
    "L4" = 65.9
    "M4" = 66.7
    "N4" = 67.5
    "O4" = 67.6
    "P4" = 67.3
    "Q4" = 67.1
    "R4" = 67.7
    "S4" = 69.2
    "T4" = 71.1
    "U4" = 71.6
    "V4" = 72.3
    "W4" = 72.9
    "X4" = 73.5
    "Y4" = 75
}

foreach ($addr in $row4.Keys) {
    $ws.Range($addr).Value = $row4[$addr]
}

# The refreshed series now runs one column further, so row 4 grows a new
# (still empty) Z cell; give it the same left/right divider border used
# elsewhere in the sheet so it reads as a trailing boundary column.
$ws.Range("Z4").Borders(7).LineStyle = 1
$ws.Range("Z4").Borders(10).LineStyle = 1

# Keep the selection where the author left it after entering the new figures.
$ws.Range("F11").Select() | Out-Null

Write-Host "row4 refreshed with mario lau's new data; Z4 boundary column added"
